# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh
# the "Latest Handoff Datetime" timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-22 16:39:46"

# zh-cn sheet: ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-22 16:39:41"

# de-de sheet: ... | Status (C) | ... | Latest Handoff Datetime (E) | ...
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-22 16:39:46"

$wb.Save()
